$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 739.619
$ws.Range("I98").Value = 627.6
$ws.Range("K98").Value = 627.6
$ws.Range("M98").Value = 870.4
$ws.Range("H112").Value = 1117.9546
$ws.Range("I112").Value = 899.5
$ws.Range("J112").Value = 1166.5
$ws.Range("K112").Value = 2698.5
$ws.Range("L112").Value = 3499.5
$ws.Range("M112").Value = -1590.5
$ws.Range("N112").Value = -5715.5
$ws.Range("H122").Value = 739.619
$ws.Range("I122").Value = 627.6
$ws.Range("K122").Value = 1882.8
$ws.Range("M122").Value = 567.1999999999998
$ws.Range("H127").Value = 1133.0
$ws.Range("I127").Value = 683.4
$ws.Range("K127").Value = 2050.2
$ws.Range("M127").Value = 2909.8
$ws.Range("H129").Value = 7620.4688
$ws.Range("I129").Value = 357.3
$ws.Range("J129").Value = 10921.909
$ws.Range("K129").Value = 1071.9
$ws.Range("L129").Value = 32765.727
$ws.Range("M129").Value = 3928.1
$ws.Range("N129").Value = -42765.727
$ws.Range("H138").Value = 2976.2952
$ws.Range("J138").Value = 3661.1
$ws.Range("L138").Value = 10983.3
$ws.Range("N138").Value = -21263.3
$ws.Range("H141").Value = 4715.6523
$ws.Range("I141").Value = 3603.158
$ws.Range("J141").Value = 10000.0
$ws.Range("K141").Value = 10809.474
$ws.Range("L141").Value = 30000.0
$ws.Range("M141").Value = -5629.474
$ws.Range("N141").Value = -40360.0

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1800.238
$ws.Range("I45").Value = 1871.909
$ws.Range("J45").Value = 1721.4
$ws.Range("K45").Value = 1871.909
$ws.Range("L45").Value = 1721.4
$ws.Range("M45").Value = -1494.909
$ws.Range("N45").Value = -2475.4
$ws.Range("H74").Value = 8336262.0
$ws.Range("I74").Value = 13891918.0
$ws.Range("J74").Value = 2778.25
$ws.Range("K74").Value = 13891918.0
$ws.Range("L74").Value = 2778.25
$ws.Range("M74").Value = -13891044.0
$ws.Range("N74").Value = -4526.25
$ws.Range("H77").Value = 8336262.0
$ws.Range("I77").Value = 13891918.0
$ws.Range("J77").Value = 2778.25
$ws.Range("K77").Value = 69459590.0
$ws.Range("L77").Value = 13891.25
$ws.Range("M77").Value = -69455222.0
$ws.Range("N77").Value = -22627.25
$ws.Range("H139").Value = 31651.166
$ws.Range("J139").Value = 31651.166
$ws.Range("L139").Value = 31651.166
$ws.Range("N139").Value = -41931.166

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 13983.333
$ws.Range("J63").Value = 13983.333
$ws.Range("L63").Value = 13983.333
$ws.Range("N63").Value = -15355.333
$ws.Range("H66").Value = 13983.333
$ws.Range("J66").Value = 13983.333
$ws.Range("L66").Value = 41949.999
$ws.Range("N66").Value = -48813.999
$ws.Range("H132").Value = 62502436.0
$ws.Range("I132").Value = 142860140.0
$ws.Range("J132").Value = 1994.1111
$ws.Range("K132").Value = 428580420.0
$ws.Range("L132").Value = 5982.3333
$ws.Range("M132").Value = -428577890.0
$ws.Range("N132").Value = -11042.3333

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 300.0
$ws.Range("I49").Value = 300.0
$ws.Range("J49").Value = 0.0
$ws.Range("K49").Value = 900.0
$ws.Range("L49").Value = 0.0
$ws.Range("N49").ClearContents()
$ws.Range("M49").Value = -744.0
$ws.Range("H80").Value = 988.6923
$ws.Range("J80").Value = 988.6923
$ws.Range("L80").Value = 2966.0769
$ws.Range("N80").Value = -4838.0769
$ws.Range("H83").Value = 988.6923
$ws.Range("J83").Value = 988.6923
$ws.Range("L83").Value = 8898.2307
$ws.Range("N83").Value = -18258.2307
$ws.Range("H113").Value = 664.2258
$ws.Range("I113").Value = 417.4375
$ws.Range("J113").Value = 927.4667
$ws.Range("K113").Value = 1252.3125
$ws.Range("L113").Value = 2782.4001
$ws.Range("M113").Value = 917.6875
$ws.Range("N113").Value = -7122.4001
$ws.Range("H124").Value = 5508.25
$ws.Range("I124").Value = 2000.0
$ws.Range("J124").Value = 9016.5
$ws.Range("K124").Value = 6000.0
$ws.Range("L124").Value = 27049.5
$ws.Range("M124").Value = -1090.0
$ws.Range("N124").Value = -36869.5
$ws.Range("H131").Value = 818.625
$ws.Range("I131").Value = 319.0
$ws.Range("J131").Value = 1118.4
$ws.Range("K131").Value = 957.0
$ws.Range("L131").Value = 3355.2
$ws.Range("M131").Value = 4083.0
$ws.Range("N131").Value = -13435.2
$ws.Range("H133").Value = 3710.0
$ws.Range("I133").Value = 3710.0
$ws.Range("J133").Value = 0.0
$ws.Range("K133").Value = 11130.0
$ws.Range("L133").Value = 0.0
$ws.Range("M133").Value = -6070.0
$ws.Range("N133").ClearContents()

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.5625
$ws.Range("I2").Value = 86.666664
$ws.Range("J2").Value = 120.71429
$ws.Range("K2").Value = 86.666664
$ws.Range("L2").Value = 120.71429
$ws.Range("M2").Value = 26.333336
$ws.Range("N2").Value = -346.71429
$ws.Range("H69").Value = 28000.0
$ws.Range("J69").Value = 28000.0
$ws.Range("L69").Value = 28000.0
$ws.Range("N69").Value = -29498.0
$ws.Range("H72").Value = 28000.0
$ws.Range("J72").Value = 28000.0
$ws.Range("L72").Value = 84000.0
$ws.Range("N72").Value = -91488.0

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2834.6667
$ws.Range("I61").Value = 2451.1
$ws.Range("J61").Value = 4752.5
$ws.Range("K61").Value = 2451.1
$ws.Range("L61").Value = 4752.5
$ws.Range("M61").Value = -2249.1
$ws.Range("N61").Value = -5156.5
$ws.Range("H113").Value = 2834.6667
$ws.Range("I113").Value = 2451.1
$ws.Range("J113").Value = 4752.5
$ws.Range("K113").Value = 2451.1
$ws.Range("L113").Value = 4752.5
$ws.Range("M113").Value = -281.0999999999999
$ws.Range("N113").Value = -9092.5
$ws.Range("H132").Value = 8069760.0
$ws.Range("I132").Value = 11911016.0
$ws.Range("K132").Value = 35733048.0
$ws.Range("M132").Value = -35730518.0
$ws.Range("H136").Value = 13340159.0
$ws.Range("I136").Value = 30314586.0
$ws.Range("J136").Value = 3108.5
$ws.Range("K136").Value = 90943758.0
$ws.Range("L136").Value = 9325.5
$ws.Range("M136").Value = -90941208.0
$ws.Range("N136").Value = -14425.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 26374.5
$ws.Range("J63").Value = 26374.5
$ws.Range("L63").Value = 26374.5
$ws.Range("N63").Value = -27622.5
$ws.Range("H66").Value = 26374.5
$ws.Range("J66").Value = 26374.5
$ws.Range("L66").Value = 79123.5
$ws.Range("N66").Value = -85363.5
$ws.Range("H132").Value = 1910.4517
$ws.Range("I132").Value = 805.8421
$ws.Range("J132").Value = 3659.4167
$ws.Range("K132").Value = 2417.5263
$ws.Range("L132").Value = 10978.2501
$ws.Range("M132").Value = 112.4737
$ws.Range("N132").Value = -16038.2501
$ws.Range("H136").Value = 6303625.5
$ws.Range("I136").Value = 9454343.0
$ws.Range("J136").Value = 2191.4285
$ws.Range("K136").Value = 28363029.0
$ws.Range("L136").Value = 6574.2855
$ws.Range("M136").Value = -28360479.0
$ws.Range("N136").Value = -11674.2855
